$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Row 24: new review for com.sugar.powerfulquotes ------------------
$ws.Range("A24").Value = "com.sugar.powerfulquotes"
$ws.Range("B24").Value = "powerful quotes"

# Column C/D carry the "email/recovery" look (font Calibri, centered) that
# the rest of the table uses - pick it up from the row above before typing
# the new values, same as Excel's "pick up formatting from adjacent cell"
# behavior when extending a table.
$ws.Range("C23").Copy()
$ws.Range("C24").PasteSpecial(-4122)
$ws.Range("C24").Value = "orenatias858@gmail.com"

$ws.Range("D23").Copy()
$ws.Range("D24").PasteSpecial(-4122)
$ws.Range("D24").Value = "dan624655@gmail.com"

$ws.Range("E24").Value = "27/5/2019 15:59"
$ws.Range("F24").Value = "great app"
$ws.Range("G24").Value = "no"

# ---- Row 25: second new review for com.sugar.powerfulquotes -----------
$ws.Range("A25").Value = "com.sugar.powerfulquotes"
$ws.Range("B25").Value = "powerful quotes"

$ws.Range("C23").Copy()
$ws.Range("C25").PasteSpecial(-4122)
$ws.Range("C25").Value = "rabuhav25@gmail.com "

$ws.Range("D23").Copy()
$ws.Range("D25").PasteSpecial(-4122)
$ws.Range("D25").Value = "itaisenior@gmail.com"

$ws.Range("E25").Value = "27/5/2019 15:59"
$ws.Range("F25").Value = "like this quotes app"
# (row 25 has no G value - left blank, matching the source edit)

# Recovery-email column on row 25 is a live mailto link, same as the other
# recovery-email cells in the sheet.
$target = $ws.Range("D25")
$ws.Hyperlinks.Add($target, "mailto:itaisenior@gmail.com", [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, "itaisenior@gmail.com")
# Hyperlinks.Add stamps the built-in "Hyperlink" style (blue/underline) on
# the cell; the sheet's other recovery-email links keep the plain look, so
# restore the shared C/D formatting after linking.
$ws.Range("D23").Copy()
$ws.Range("D25").PasteSpecial(-4122)

# Clear the copy marquee and leave the selection where data entry stopped.
$excel.CutCopyMode = 0
$ws.Range("F26").Select()
